$d = $word.ActiveDocument

# ===========================================================================
# Part 1: merge split runs that should become a single run within a
# paragraph (no visual change, just simplifying run structure / small
# wording tweaks). In every case the surviving run is deliberately chosen
# to be one whose *original* text had no leading/trailing whitespace, so
# that re-inserting the removed text via InsertBefore/InsertAfter does not
# leave a stray xml:space="preserve" attribute behind.
# ===========================================================================

function Merge-TwoRuns($paraIndex, $firstRunLen, $prefixText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pStart = $r.Start
    $delRng = $d.Range($pStart, $pStart + $firstRunLen)
    $delRng.Delete()
    $keepRng = $d.Range($pStart, $pStart)
    $keepRng.InsertBefore($prefixText)
}

function Merge-ThreeRuns-KeepMiddle($paraIndex, $run1Len, $run2Len, $prefixText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pStart = $r.Start
    $pEnd = $r.End - 1

    $keepStart = $pStart + $run1Len
    $keepEnd = $keepStart + $run2Len

    $suffixRng = $d.Range($keepEnd, $pEnd)
    $suffixText = $suffixRng.Text
    $suffixRng.Delete()

    $prefixRng = $d.Range($pStart, $keepStart)
    $prefixRng.Delete()

    $keepRng = $d.Range($pStart, $pStart + $run2Len)
    $keepRng.InsertBefore($prefixText)

    $afterPos = $pStart + $prefixText.Length + $run2Len
    $afterRng = $d.Range($afterPos, $afterPos)
    $afterRng.InsertAfter($suffixText)
}

function Merge-ThreeRuns-KeepLast($paraIndex, $deleteLen, $prefixText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pStart = $r.Start
    $delRng = $d.Range($pStart, $pStart + $deleteLen)
    $delRng.Delete()
    $keepRng = $d.Range($pStart, $pStart)
    $keepRng.InsertBefore($prefixText)
}

# ---------------------------------------------------------------------------
# 1. "- столбцовые: " + "C-Store;"  ->  "- столбцовые: C-Store;"
# ---------------------------------------------------------------------------
Merge-TwoRuns 11 "- столбцовые: ".Length "- столбцовые: "

# ---------------------------------------------------------------------------
# 2. "- " + "семейство столбцов:" + " Cassandra, Hbase, Hypertable;"
#    -> "- семейство столбцов: Cassandra, Hbase, Hypertable;"
# ---------------------------------------------------------------------------
Merge-ThreeRuns-KeepMiddle 12 "- ".Length "семейство столбцов:".Length "- "

# ---------------------------------------------------------------------------
# 3. "Столбцовые (также известны как «семейства столбцов», англ. column-wide) "
#    + "СУБД возникли ... [Фаулер]."
#    -> "Столбцовые СУБД возникли ... [Фаулер]."
# ---------------------------------------------------------------------------
Merge-TwoRuns 21 "Столбцовые (также известны как «семейства столбцов», англ. column-wide) ".Length "Столбцовые "

# ---------------------------------------------------------------------------
# 4. "Отличительной особенностью б" + "аз данных «семейств столбцов» (англ. column-family, тж. wide-column) "
#    + "является распределение ... [Редмонд]"
#    -> single merged run (text unchanged).
# ---------------------------------------------------------------------------
$run23_1 = "Отличительной особенностью б".Length
$run23_2 = "аз данных «семейств столбцов» (англ. column-family, тж. wide-column) ".Length
$deleteLen23 = $run23_1 + $run23_2
Merge-ThreeRuns-KeepLast 23 $deleteLen23 "Отличительной особенностью баз данных «семейств столбцов» (англ. column-family, тж. wide-column) "

# ---------------------------------------------------------------------------
# 5. "Хранилища " + "данных" + " «ключ-значение» состоят ... [Hoffner]"
#    -> single merged run (text unchanged).
# ---------------------------------------------------------------------------
Merge-ThreeRuns-KeepMiddle 24 "Хранилища ".Length "данных".Length "Хранилища "

Write-Host "Merges complete"

# ===========================================================================
# Part 2: add the two new paragraphs about graph databases after the
# paragraph ending "...для многих веб-разработчиков [Harrison G. Next Gen.]."
# (currently paragraph 26).
#
# To get multiple separate <w:r> runs (matching the target XML) instead of
# one run, the whole paragraph text is typed in one go and then each
# internal boundary is "split" by selecting the *exact, non-collapsed*
# range of text that should become its own run and toggling Bold on/off on
# it. (Toggling formatting on a zero-length/collapsed range was found to
# corrupt neighbouring paragraphs' formatting in this runtime, so a
# concrete, fully-bounded range must always be used.)
# ===========================================================================

function Split-Run($startPos, $endPos) {
    $sr = $d.Range($startPos, $endPos)
    $sr.Bold = 1
    $sr.Bold = 0
}

# --- New paragraph A: "Графовая база данных ... [Фаулер NoSQL]." ---------
$pAnchor = $d.Paragraphs.Item(26)
$rAnchor = $pAnchor.Range
$rAnchor.InsertParagraphAfter()

$paraA = $d.Paragraphs.Item(27)
$rngA = $paraA.Range
$rngA.End = $rngA.End - 1
$startA = $rngA.Start

$a1 = "Графовая база данных состоит из набора вершин (узлов, сущностей) и граней (связей, отношений, рёбер). Узлы воспринимаются как объекты со свойствами, между которыми моделируются отношения с помощью граней, которые также могут иметь свойства. "
$a2 = "Отношения имеют направления, на их основе происходит организация узлов, что позволяет единожды записать данные и затем по-разному их интерпретировать"
$a3 = " [Jordan] "
$a4 = "[Фаулер NoSQL]."
$fullA = $a1 + $a2 + $a3 + $a4
$rngA.InsertAfter($fullA)

$aPos2Start = $startA + $a1.Length
$aPos2End = $aPos2Start + $a2.Length
$aPos3End = $aPos2End + $a3.Length
$aPos4End = $aPos3End + $a4.Length

Split-Run $aPos2Start $aPos2End
Split-Run $aPos2End $aPos3End
Split-Run $aPos3End $aPos4End

Write-Host "Paragraph A (graph intro) inserted"

# --- New paragraph B: "Графовая структура ... [Bruggen]." -----------------
$pAnchorB = $d.Paragraphs.Item(27)
$rAnchorB = $pAnchorB.Range
$rAnchorB.InsertParagraphAfter()

$paraB = $d.Paragraphs.Item(28)
$rngB = $paraB.Range
$rngB.End = $rngB.End - 1
$startB = $rngB.Start

$b1 = "Графовая структура позволяет представить данные в более естественном виде без искажений, как это может произойти в реляционных базах данных, а также применить различные типы графовых алгоритмов к этим данным. "
$b2 = "Одна из ключевых особенностей графовых БД — возможность обхода графа по его узлам и граням, перемещения от одного узла к другому, следуя направленным отношениям. Эта возможность называется «index free adjacency» (примерно переводится как смежность без индекса), смысл которой заключается в поиске прилежащих узлов без использования поиска по индексу, что значительным образом сказывается на производительности [Bruggen]."
$fullB = $b1 + $b2
$rngB.InsertAfter($fullB)

$bPos2Start = $startB + $b1.Length
$bPos2End = $bPos2Start + $b2.Length

Split-Run $bPos2Start $bPos2End

Write-Host "Paragraph B (graph structure) inserted"
